$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Keywords")

# Insert a new data row into the table, right above the "Resid" row (old row 103),
# shifting all subsequent rows (and the table's last-row formula binding) down by one.
$ws.Rows.Item(103).Insert()

# Grow the table definition so it covers the newly inserted row.
$lo = $ws.ListObjects.Item(1)
$lo.Resize($ws.Range("A1:C145"))

$ws.Range("A103").Value = "Novot"
$ws.Range("B103").Formula = "=LEN(Cluster_Keywords[[#This Row],[Stem]])"
$ws.Range("C103").Value = "Housing"

# Match the style of the neighbouring table rows (style index 3 in the original file).
$ws.Range("A103:C103").Style = $ws.Range("A104:C104").Style

# Inserting a row in the middle of the table leaves the calculated-column formula
# in the table's last row pointing at a stale structured reference; re-enter it so
# it recalculates correctly against the new table extent.
$ws.Range("B145").Formula = "=LEN(Cluster_Keywords[[#This Row],[Stem]])"
